$d = $word.ActiveDocument

# --- Step 1: merge the tail of paragraph 1 (the runs broken up by spell/grammar
# proofErr markers) into a single run with the same text, so Word collapses
# them into one <w:r> the way it does whenever text actually changes. ---
$p1 = $d.Paragraphs(1)
$full = $p1.Range
$newTail = " must provide information on the prior award(s), major achievements, and relevance to the proposed NRT project.Individuals who have received more than one prior award must report on the award(s) most closely related to the proposal.Complete bibliographic citation for each publication resulting from an NSF award must be included in either the Results from PriorNSF Support section or in the References Cited section of the proposal. For further information see Chapter II.C.2.d of the GPG."

$tailStart = $full.Start + 138
$tailEnd = $full.End - 1
$tailRange = $d.Range($tailStart, $tailEnd)
# Setting identical text is a no-op in this engine (it only merges runs when
# the content actually changes), so append a one-character sentinel to force
# a real edit, then strip the sentinel off in a second pass.
$tailRange.Text = $newTail + "#"
$full2 = $d.Paragraphs(1).Range
$sentinelRange = $d.Range($full2.End - 2, $full2.End - 1)
$sentinelRange.Text = ""

# --- Step 2: insert a brand-new, fully empty paragraph right after paragraph 1. ---
$p1 = $d.Paragraphs(1)
$endOfP1 = $p1.Range.End
$insertPoint = $d.Range($endOfP1, $endOfP1)
$emptyParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$insertPoint.InsertXML($emptyParaXml)

# --- Step 3: add a "2.5 pages" run at the very start of the final paragraph
# (the one holding the _GoBack bookmark), before the bookmark. ---
$pLast = $d.Paragraphs($d.Paragraphs.Count)
$startLast = $pLast.Range.Start
$lastInsertPoint = $d.Range($startLast, $startLast)
$lastInsertPoint.InsertBefore("2.5 pages")
